$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# xlPasteFormats = -4122, xlPasteValues = -4123

# --- Row 12: repeat of the header row (row 1), columns E:L only ---
$ws.Range("E1:L1").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("E1:L1").Copy()
$ws.Range("E12").PasteSpecial(-4123)
$ws.Rows("12:12").RowHeight = 52.5

# --- Row 13: clone of row 2's format (D:L), then patch in the extra cells ---
$ws.Range("D2:L2").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Rows("13:13").RowHeight = 52.5
$ws.Range("D13").Value = 1
$ws.Range("F13").Value = "C"
$ws.Range("G13").Value = "A"
$ws.Range("H13").Value = "B"
$ws.Range("I13").Value = "B"
$ws.Range("J13").Value = "A"
$ws.Range("K13").Value = "C"

# --- Row 14: exact clone of row 3 (D:L) ---
$ws.Range("D3:L3").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D3:L3").Copy()
$ws.Range("D14").PasteSpecial(-4123)
$ws.Rows("14:14").RowHeight = 52.5

# --- Row 15: clone of row 4's format (D:L), then patch in E/L ---
$ws.Range("D4:L4").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Rows("15:15").RowHeight = 52.5
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = "A"
$ws.Range("L15").Value = "A"

# --- Row 16: clone of row 5's format (D:L), then patch in E/L ---
$ws.Range("D5:L5").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Rows("16:16").RowHeight = 52.5
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = "B"
$ws.Range("H16").Value = "W"
$ws.Range("I16").Value = "S"
$ws.Range("L16").Value = "B"

# --- Row 17: clone of row 6's format (D:L), then patch in E/L ---
$ws.Range("D6:L6").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Rows("17:17").RowHeight = 52.5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = "B"
$ws.Range("H17").Value = "S"
$ws.Range("I17").Value = "W"
$ws.Range("L17").Value = "B"

# --- Row 18: clone of row 7's format (D:L), then patch in E/L ---
$ws.Range("D7:L7").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Rows("18:18").RowHeight = 52.5
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = "A"
$ws.Range("L18").Value = "A"

# --- Row 19: exact clone of row 8 (D:L) ---
$ws.Range("D8:L8").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D8:L8").Copy()
$ws.Range("D19").PasteSpecial(-4123)
$ws.Rows("19:19").RowHeight = 52.5

# --- Row 20: clone of row 9's format (D:L), then patch in the extra cells ---
$ws.Range("D9:L9").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Rows("20:20").RowHeight = 52.5
$ws.Range("D20").Value = 8
$ws.Range("F20").Value = "C"
$ws.Range("G20").Value = "A"
$ws.Range("H20").Value = "B"
$ws.Range("I20").Value = "B"
$ws.Range("J20").Value = "A"
$ws.Range("K20").Value = "C"

$excel.CutCopyMode = $false

# --- View state: match the committed sheet view (active cell S20) ---
$ws.Activate()
$ws.Range("S20").Select() | Out-Null

# --- New, empty second sheet, placed after Tabelle1 ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Tabelle2"

# Keep Tabelle1 as the active/selected sheet
$ws1.Activate()
